$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.043483244119445
$ws.Range("D2").Value = 1.051866369584021
$ws.Range("E2").Value = 1.041476039531331
$ws.Range("F2").Value = 1.062228763669788
$ws.Range("I2").Value = 1.044373970303589
$ws.Range("J2").Value = 1.048553449909775
$ws.Range("K2").Value = 1.054616684860785
$ws.Range("L2").Value = 1.04425544724651
$ws.Range("M2").Value = 1.064950678572257
$ws.Range("N2").Value = 1.020019831758416
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.044477228607705
$ws.Range("D3").Value = 1.052671590926876
$ws.Range("E3").Value = 1.04232363705239
$ws.Range("F3").Value = 1.063191708835326
$ws.Range("I3").Value = 1.044655656931538
$ws.Range("J3").Value = 1.049194079399179
$ws.Range("K3").Value = 1.05523465720853
$ws.Range("L3").Value = 1.044913583143965
$ws.Range("M3").Value = 1.065728026206772
$ws.Range("N3").Value = 1.020237273422024
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.045120796687535
$ws.Range("D4").Value = 1.053192975547602
$ws.Range("E4").Value = 1.042872765755254
$ws.Range("F4").Value = 1.063815519601944
$ws.Range("I4").Value = 1.044836973774907
$ws.Range("J4").Value = 1.049608392617939
$ws.Range("K4").Value = 1.055634224643341
$ws.Range("L4").Value = 1.045339473285014
$ws.Range("M4").Value = 1.06623111874233
$ws.Range("N4").Value = 1.02037778952169
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.045391446325477
$ws.Range("D5").Value = 1.053412248967678
$ws.Range("E5").Value = 1.04310378031729
$ws.Range("F5").Value = 1.064077941132445
$ws.Range("I5").Value = 1.044912970730451
$ws.Range("J5").Value = 1.04978251722917
$ws.Range("K5").Value = 1.055802129701518
$ws.Range("L5").Value = 1.045518524388832
$ws.Range("M5").Value = 1.066442641045123
$ws.Range("N5").Value = 1.02043681828751
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.045436895062367
$ws.Range("D6").Value = 1.053449070794646
$ws.Range("E6").Value = 1.043142578079948
$ws.Range("F6").Value = 1.064122012886433
$ws.Range("I6").Value = 1.044925717533148
$ws.Range("J6").Value = 1.049811750407771
$ws.Range("K6").Value = 1.055830317390495
$ws.Range("L6").Value = 1.045548588237121
$ws.Range("M6").Value = 1.066478157845283
$ws.Range("N6").Value = 1.020446726874958
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.045124412754157
$ws.Range("D7").Value = 1.05319590516371
$ws.Range("E7").Value = 1.042875851952441
$ws.Range("F7").Value = 1.06381902541882
$ws.Range("I7").Value = 1.044837990148908
$ws.Range("J7").Value = 1.049610719487378
$ws.Range("K7").Value = 1.055636468487132
$ws.Range("L7").Value = 1.045341865748585
$ws.Range("M7").Value = 1.066233945027622
$ws.Range("N7").Value = 1.020378578441064
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.043819083946909
$ws.Range("D8").Value = 1.052138424255614
$ws.Range("E8").Value = 1.041762348371233
$ws.Range("F8").Value = 1.062554045583421
$ws.Range("I8").Value = 1.044469364729339
$ws.Range("J8").Value = 1.048769998109898
$ws.Range("K8").Value = 1.054825593562275
$ws.Range("L8").Value = 1.044477860191154
$ws.Range("M8").Value = 1.065213366273413
$ws.Range("N8").Value = 1.020093354911397
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.04152197015955
$ws.Range("D9").Value = 1.050277755666849
$ws.Range("E9").Value = 1.039805441206605
$ws.Range("F9").Value = 1.060330558454541
$ws.Range("I9").Value = 1.043812518296459
$ws.Range("J9").Value = 1.047286911065883
$ws.Range("K9").Value = 1.053394450371197
$ws.Range("L9").Value = 1.042955654913875
$ws.Range("M9").Value = 1.063415757438564
$ws.Range("N9").Value = 1.019589366736693
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.039992641678471
$ws.Range("D10").Value = 1.049039220660363
$ws.Range("E10").Value = 1.038504415107861
$ws.Range("F10").Value = 1.058852041546026
$ws.Range("I10").Value = 1.04336975061868
$ws.Range("J10").Value = 1.046297134349964
$ws.Range("K10").Value = 1.052438869504675
$ws.Range("L10").Value = 1.041941089501485
$ws.Range("M10").Value = 1.062217931319337
$ws.Range("N10").Value = 1.019252461637645
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.039330925903068
$ws.Range("D11").Value = 1.048503387522468
$ws.Range("E11").Value = 1.037941918297257
$ws.Range("F11").Value = 1.058212743578758
$ws.Range("I11").Value = 1.043176877391684
$ws.Range("J11").Value = 1.045868309282419
$ws.Range("K11").Value = 1.052024749534467
$ws.Range("L11").Value = 1.041501837735689
$ws.Range("M11").Value = 1.061699408016274
$ws.Range("N11").Value = 1.019106365597927
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.039085209825089
$ws.Range("D12").Value = 1.048304425457984
$ws.Range("E12").Value = 1.037733111566584
$ws.Range("F12").Value = 1.05797541750915
$ws.Range("I12").Value = 1.043105063006309
$ws.Range("J12").Value = 1.045708988369803
$ws.Range("K12").Value = 1.05187087549545
$ws.Range("L12").Value = 1.041338690096383
$ws.Range("M12").Value = 1.061506827752576
$ws.Range("N12").Value = 1.019052067203059
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.039137913373364
$ws.Range("D13").Value = 1.048347100302384
$ws.Range("E13").Value = 1.0377778954368
$ws.Range("F13").Value = 1.058026318511184
$ws.Range("I13").Value = 1.043120475243827
$ws.Range("J13").Value = 1.045743164873522
$ws.Range("K13").Value = 1.051903884315625
$ws.Range("L13").Value = 1.041373685353212
$ws.Range("M13").Value = 1.061548135846568
$ws.Range("N13").Value = 1.019063715827952
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.039310613410342
$ws.Range("D14").Value = 1.048486939817835
$ws.Range("E14").Value = 1.037924655617478
$ws.Range("F14").Value = 1.05819312331778
$ws.Range("I14").Value = 1.043170944712627
$ws.Range("J14").Value = 1.045855140500545
$ws.Range("K14").Value = 1.052012031303703
$ws.Range("L14").Value = 1.041488351692114
$ws.Range("M14").Value = 1.061683488810176
$ws.Range("N14").Value = 1.019101877922971
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.03941702950562
$ws.Range("D15").Value = 1.048573108893692
$ws.Range("E15").Value = 1.038015096616083
$ws.Range("F15").Value = 1.058295915535697
$ws.Range("I15").Value = 1.043202017747512
$ws.Range("J15").Value = 1.045924127604169
$ws.Range("K15").Value = 1.052078657436951
$ws.Range("L15").Value = 1.041559002758648
$ws.Range("M15").Value = 1.061766887229014
$ws.Range("N15").Value = 1.019125386642132
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.040036568013983
$ws.Range("D16").Value = 1.049074791933196
$ws.Range("E16").Value = 1.038541764291719
$ws.Range("F16").Value = 1.058894488896504
$ws.Range("I16").Value = 1.04338252674259
$ws.Range("J16").Value = 1.046325588946535
$ws.Range("K16").Value = 1.052466346029896
$ws.Range("L16").Value = 1.04197024255928
$ws.Range("M16").Value = 1.062252347094915
$ws.Range("N16").Value = 1.019262153070853
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.040425320622003
$ws.Range("D17").Value = 1.049389608500678
$ws.Range("E17").Value = 1.038872358921712
$ws.Range("F17").Value = 1.059270202386393
$ws.Range("I17").Value = 1.043495447136254
$ws.Range("J17").Value = 1.046577349847885
$ws.Range("K17").Value = 1.052709440401757
$ws.Range("L17").Value = 1.04222821948087
$ws.Range("M17").Value = 1.062556902208949
$ws.Range("N17").Value = 1.019347885992276
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.040652121029322
$ws.Range("D18").Value = 1.049573280060818
$ws.Range("E18").Value = 1.039065271706981
$ws.Range("F18").Value = 1.05948943740543
$ws.Range("I18").Value = 1.043561200524135
$ws.Range("J18").Value = 1.046724174018328
$ws.Range("K18").Value = 1.052851199693308
$ws.Range("L18").Value = 1.042378698962709
$ws.Range("M18").Value = 1.062734557824755
$ws.Range("N18").Value = 1.01939787189454
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.040729462164631
$ws.Range("D19").Value = 1.049635914813765
$ws.Range("E19").Value = 1.039131063949257
$ws.Range("F19").Value = 1.059564205725876
$ws.Range("I19").Value = 1.0435836018628
$ws.Range("J19").Value = 1.046774233225783
$ws.Range("K19").Value = 1.052899530237288
$ws.Range("L19").Value = 1.042430009556706
$ws.Range("M19").Value = 1.062795136108677
$ws.Range("N19").Value = 1.019414912281484
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.040383606220524
$ws.Range("D20").Value = 1.049355827064376
$ws.Range("E20").Value = 1.038836880702724
$ws.Range("F20").Value = 1.059229882806466
$ws.Range("I20").Value = 1.043483343342483
$ws.Range("J20").Value = 1.046550340732868
$ws.Range("K20").Value = 1.052683362140588
$ws.Range("L20").Value = 1.042200540382127
$ws.Range("M20").Value = 1.062524224896995
$ws.Range("N20").Value = 1.019338689795871
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.039259755490067
$ws.Range("D21").Value = 1.048445758619408
$ws.Range("E21").Value = 1.037881434815597
$ws.Range("F21").Value = 1.058143999651586
$ws.Range("I21").Value = 1.043156087475772
$ws.Range("J21").Value = 1.045822167463051
$ws.Range("K21").Value = 1.051980186125566
$ws.Range("L21").Value = 1.041454585030569
$ws.Range("M21").Value = 1.061643630111406
$ws.Range("N21").Value = 1.019090641013019
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.038553577920274
$ws.Range("D22").Value = 1.047873969243588
$ws.Range("E22").Value = 1.037281458566821
$ws.Range("F22").Value = 1.057462058770831
$ws.Range("I22").Value = 1.042949329833282
$ws.Range("J22").Value = 1.045364126973178
$ws.Range("K22").Value = 1.051537774174649
$ws.Range("L22").Value = 1.040985631662064
$ws.Range("M22").Value = 1.061090095114705
$ws.Range("N22").Value = 1.018934498949842
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.038927894903894
$ws.Range("D23").Value = 1.048177046716038
$ws.Range("E23").Value = 1.037599445850299
$ws.Range("F23").Value = 1.057823492656955
$ws.Range("I23").Value = 1.04305903052752
$ws.Range("J23").Value = 1.045606962500695
$ws.Range("K23").Value = 1.051772333061012
$ws.Range("L23").Value = 1.041234226933222
$ws.Range("M23").Value = 1.061383521893332
$ws.Range("N23").Value = 1.019017290142248
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.040402455015907
$ws.Range("D24").Value = 1.049371091304278
$ws.Range("E24").Value = 1.038852911527815
$ws.Range("F24").Value = 1.059248101218195
$ws.Range("I24").Value = 1.043488812869389
$ws.Range("J24").Value = 1.046562545062618
$ws.Range("K24").Value = 1.052695145887964
$ws.Range("L24").Value = 1.042213047357205
$ws.Range("M24").Value = 1.062538990325189
$ws.Range("N24").Value = 1.019342845225041
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.042115464993377
$ws.Range("D25").Value = 1.05075845095103
$ws.Range("E25").Value = 1.040310722493766
$ws.Range("F25").Value = 1.060904716680806
$ws.Range("I25").Value = 1.043983189102252
$ws.Range("J25").Value = 1.047670513193927
$ws.Range("K25").Value = 1.053764700392008
$ws.Range("L25").Value = 1.043349143009977
$ws.Range("M25").Value = 1.063880384218365
$ws.Range("N25").Value = 1.01971982215533
